# Updated vehicle design sheet to include brake parameters
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New text cells (columns A and C) -------------------------------------
# Written in this specific order so new shared-string table entries are
# appended in the same sequence the original authoring session produced.
$ws.Range("A71").Value = "piston_radius_front"
$ws.Range("A72").Value = "piston_radius_rear"
$ws.Range("A74").Value = "num_pistons_rear"
$ws.Range("A75").Value = "pad_friction_front"
$ws.Range("A76").Value = "pad_friction_rear"
$ws.Range("A77").Value = "max_pedal_force"
$ws.Range("C73").Value = "unitless"
$ws.Range("C74").Value = "unitless"
$ws.Range("C75").Value = "unitless"
$ws.Range("C76").Value = "unitless"
$ws.Range("A78").Value = "disc_radius_front"
$ws.Range("A79").Value = "disc_radius_rear"
$ws.Range("A80").Value = "pad_height_front"
$ws.Range("A81").Value = "pad_height_rear"
$ws.Range("A73").Value = "num_pistons_front"
$ws.Range("C71").Value = "m, 0.5 in"
$ws.Range("C72").Value = "m, 0.5 in"
$ws.Range("C77").Value = "N, 150 lbf"
$ws.Range("A82").Value = "mc_diameter_front"
$ws.Range("A83").Value = "mc_diameter_rear"
$ws.Range("A84").Value = "balance_bar_ratio_front"
$ws.Range("C82").Value = "m, 0.625 in"
$ws.Range("C83").Value = "m, 0.9375 in"
$ws.Range("C80").Value = "m, 1 in"
$ws.Range("C81").Value = "m, 1 in"
$ws.Range("C78").Value = "m, 5.125 in"
$ws.Range("C79").Value = "m, 5.125 in"
$ws.Range("C84").Value = "0 to 1, from brake design spreadsheet"

# --- New numeric values (column B), styled "Neutral" like the other --------
# "value needs confirmed" rows on this sheet.
$ws.Range("B71").Value = 0.013
$ws.Range("B72").Value = 0.013
$ws.Range("B73").Value = 2
$ws.Range("B74").Value = 2
$ws.Range("B75").Value = 0.5
$ws.Range("B76").Value = 0.5
$ws.Range("B77").Value = 667.23
$ws.Range("B78").Value = 0.13
$ws.Range("B79").Value = 0.13
$ws.Range("B80").Value = 0.025
$ws.Range("B81").Value = 0.025
$ws.Range("B82").Value = 0.016
$ws.Range("B83").Value = 0.024
$ws.Range("B84").Value = 0.5152

$ws.Range("B71:B84").Style = "Neutral"

# --- View state: scroll position / selection -------------------------------
$ws.Application.ActiveWindow.ScrollRow = 45
$ws.Range("D78").Select()

# --- Page setup --------------------------------------------------------
$ws.PageSetup.Orientation = 1
$ws.PageSetup.HorizontalDpi = 0
$ws.PageSetup.VerticalDpi = 0

# --- Workbook window geometry -----------------------------------------
$excel.ActiveWindow.Left = 60
$excel.ActiveWindow.Width = 14240
